# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (strikeout) values for column G, keyed by row number.
$kValues = @{
    2  = 1
    3  = 0
    4  = 2
    5  = 1
    6  = 1
    7  = 0
    9  = 2
    10 = 0
    11 = 2
    12 = 0
    13 = 1
    14 = 2
    15 = 1
    16 = 0
    17 = 3
    18 = 2
    19 = 2
    20 = 2
    21 = 1
    23 = 1
    24 = 0
    27 = 2
    28 = 2
    29 = 2
    30 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
